$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

# Update the repayment strategy value in B17 to reflect the new scenario
$ws.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

# Reflect the new active selection (matches the author's saved cursor position)
$ws.Range("B17").Select()
